# Markdown inline formatting: replace the heading-styled title with plain
# body text and drop the decorative "Heading 1" styling (bold font +
# thick bottom border) that used to set the title apart, along with the
# now-empty spacer row that existed only to host the border's "thick top".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text shown in A1 and strip any special formatting back to the
# workbook's default ("Normal") cell style.
$ws.Range("A1").Value = "lorem ipsum"
$ws.Range("A1").Style = "Normal"

# Row 2 was only present to carry the heading's "thick top" border edge;
# with the heading gone it is removed entirely.
$ws.Rows("2").Delete()

# Row 1 no longer needs the taller "heading" row height / thick bottom
# border, so let it size itself back to the default again.
$ws.Rows("1").AutoFit()

# Remove the now-unused built-in "Heading 1" cell style definition.
$wb.Styles("Heading 1").Delete()

# Move the selection down to A2, matching where the cursor ends up after
# the heading row is gone.
$ws.Range("A2").Select() | Out-Null
